# no-op test
$p = $ppt.ActivePresentation
Write-Host "Slides count: $($p.Slides.Count)"
